$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates derived from the authoritative diff.
# Numeric-looking strings are forced to remain text (matching the
# original inlineStr cell type) by temporarily setting a text number
# format, then resetting the style back to Normal so no stray
# formatting is left behind on the cell.

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '311.88'
Set-TextValue $ws.Range('E2') '1.96%'
Set-TextValue $ws.Range('G2') '22'
Set-TextValue $ws.Range('D3') '37.36'
Set-TextValue $ws.Range('E3') '0.69%'
Set-TextValue $ws.Range('G3') '22'
Set-TextValue $ws.Range('D4') '5.135'
Set-TextValue $ws.Range('E4') '0.94%'
Set-TextValue $ws.Range('G4') '22'
Set-TextValue $ws.Range('D5') '0.07831'
Set-TextValue $ws.Range('E5') '1.60%'
Set-TextValue $ws.Range('G5') '22'
Set-TextValue $ws.Range('D6') '4.431'
Set-TextValue $ws.Range('E6') '1.79%'
Set-TextValue $ws.Range('G6') '22'
Set-TextValue $ws.Range('B7') 'FTXToken'
Set-TextValue $ws.Range('C7') 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range('D7') '1.909'
Set-TextValue $ws.Range('E7') '0.98%'
Set-TextValue $ws.Range('G7') '22'
Set-TextValue $ws.Range('B8') 'KuCoinToken'
Set-TextValue $ws.Range('C8') 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue $ws.Range('D8') '8.281'
Set-TextValue $ws.Range('E8') '1.07%'
Set-TextValue $ws.Range('G8') '22'
Set-TextValue $ws.Range('D9') '2.835'
Set-TextValue $ws.Range('E9') '-5.93%'
Set-TextValue $ws.Range('G9') '22'
Set-TextValue $ws.Range('D10') '0.9189'
Set-TextValue $ws.Range('E10') '0.20%'
Set-TextValue $ws.Range('G10') '22'
Set-TextValue $ws.Range('D11') '0.1197'
Set-TextValue $ws.Range('E11') '1.13%'
Set-TextValue $ws.Range('G11') '22'
Set-TextValue $ws.Range('D12') '0.1916'
Set-TextValue $ws.Range('G12') '22'
Set-TextValue $ws.Range('D13') '0.09032'
Set-TextValue $ws.Range('E13') '3.98%'
Set-TextValue $ws.Range('G13') '22'
Set-TextValue $ws.Range('D14') '0.03353'
Set-TextValue $ws.Range('E14') '-1.14%'
Set-TextValue $ws.Range('G14') '22'
Set-TextValue $ws.Range('D15') '0.09591'
Set-TextValue $ws.Range('E15') '-1.09%'
Set-TextValue $ws.Range('G15') '22'
Set-TextValue $ws.Range('D16') '0.001382'
Set-TextValue $ws.Range('E16') '1.30%'
Set-TextValue $ws.Range('G16') '22'
Set-TextValue $ws.Range('D17') '0.005698'
Set-TextValue $ws.Range('E17') '-4.47%'
Set-TextValue $ws.Range('G17') '22'
Set-TextValue $ws.Range('D18') '3.528'
Set-TextValue $ws.Range('E18') '-2.45%'
Set-TextValue $ws.Range('G18') '22'
Set-TextValue $ws.Range('D19') '0.3441'
Set-TextValue $ws.Range('E19') '0.90%'
Set-TextValue $ws.Range('G19') '22'
Set-TextValue $ws.Range('D20') '5.261'
Set-TextValue $ws.Range('E20') '4.88%'
Set-TextValue $ws.Range('G20') '22'
Set-TextValue $ws.Range('D21') '0.1283'
Set-TextValue $ws.Range('E21') '0.70%'
Set-TextValue $ws.Range('G21') '22'
Set-TextValue $ws.Range('D22') '0.2595'
Set-TextValue $ws.Range('E22') '-0.04%'
Set-TextValue $ws.Range('G22') '22'
Set-TextValue $ws.Range('D23') '0.04363'
Set-TextValue $ws.Range('E23') '0.98%'
Set-TextValue $ws.Range('G23') '22'
Set-TextValue $ws.Range('D24') '0.001251'
Set-TextValue $ws.Range('E24') '3.08%'
Set-TextValue $ws.Range('G24') '22'
Set-TextValue $ws.Range('D25') '0.004677'
Set-TextValue $ws.Range('E25') '11.08%'
Set-TextValue $ws.Range('G25') '22'
Set-TextValue $ws.Range('E26') '0.75%'
Set-TextValue $ws.Range('G26') '22'
Set-TextValue $ws.Range('D27') '0.0003995'
Set-TextValue $ws.Range('E27') '-98.11%'
Set-TextValue $ws.Range('G27') '22'
Set-TextValue $ws.Range('G28') '22'
Set-TextValue $ws.Range('G29') '22'
Set-TextValue $ws.Range('G30') '22'
Set-TextValue $ws.Range('G31') '22'
Set-TextValue $ws.Range('G32') '22'
Set-TextValue $ws.Range('G33') '22'
Set-TextValue $ws.Range('G34') '22'
Set-TextValue $ws.Range('G35') '22'
Set-TextValue $ws.Range('G36') '22'
Set-TextValue $ws.Range('G37') '22'
Set-TextValue $ws.Range('G38') '22'
Set-TextValue $ws.Range('D39') '0.02270'
Set-TextValue $ws.Range('E39') '3.27%'
Set-TextValue $ws.Range('G39') '22'
Set-TextValue $ws.Range('D40') '0.05059'
Set-TextValue $ws.Range('E40') '3.07%'
Set-TextValue $ws.Range('G40') '22'
Set-TextValue $ws.Range('D41') '0.007486'
Set-TextValue $ws.Range('E41') '-0.82%'
Set-TextValue $ws.Range('G41') '22'
Set-TextValue $ws.Range('D42') '0.009066'
Set-TextValue $ws.Range('E42') '-8.52%'
Set-TextValue $ws.Range('G42') '22'
Set-TextValue $ws.Range('D43') '0.1351'
Set-TextValue $ws.Range('E43') '1.54%'
Set-TextValue $ws.Range('G43') '22'
Set-TextValue $ws.Range('D44') '0.002004'
Set-TextValue $ws.Range('E44') '-2.91%'
Set-TextValue $ws.Range('G44') '22'
Set-TextValue $ws.Range('D45') '0.009297'
Set-TextValue $ws.Range('E45') '9.40%'
Set-TextValue $ws.Range('G45') '22'
Set-TextValue $ws.Range('D46') '0.00006625'
Set-TextValue $ws.Range('E46') '1.08%'
Set-TextValue $ws.Range('G46') '22'
Set-TextValue $ws.Range('E47') '0.01%'
Set-TextValue $ws.Range('G47') '22'
Set-TextValue $ws.Range('D48') '0.003301'
Set-TextValue $ws.Range('E48') '9.97%'
Set-TextValue $ws.Range('G48') '22'
Set-TextValue $ws.Range('D49') '0.001001'
Set-TextValue $ws.Range('E49') '-23.12%'
Set-TextValue $ws.Range('G49') '22'
Set-TextValue $ws.Range('D50') '0.00002104'
Set-TextValue $ws.Range('E50') '0.01%'
Set-TextValue $ws.Range('G50') '22'
Set-TextValue $ws.Range('D51') '0.0002004'
Set-TextValue $ws.Range('E51') '0.01%'
Set-TextValue $ws.Range('G51') '22'
